$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Latest Period" (N column) date values, preserving/adjusting highlight style ---
# Stable template cells (untouched by this edit) used to copy exact cell formatting (fill) via PasteSpecial:
#   N10 -> style 49 (yellow highlight, no border)
#   N39 -> style 48 (no highlight)

$ws.Range("N10").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 46023

$ws.Range("N39").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 46023

$ws.Range("N10").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = 46023

$ws.Range("N29").Value = 46063

$ws.Range("N30").Value = 46063

$ws.Range("N10").Copy()
$ws.Range("N33").PasteSpecial(-4122)
$ws.Range("N33").Value = 46023

$ws.Range("N47").Value = 46062

$ws.Range("N48").Value = 46062

$ws.Range("N49").Value = 46062

$ws.Range("N50").Value = 46062

$ws.Range("N52").Value = 46062

$excel.CutCopyMode = $false

# --- Update Present/Lag1..Lag4 (Q-U) numeric values ---
$ws.Range("Q3").Value = 130
$ws.Range("R3").Value = 48
$ws.Range("S3").Value = 41
$ws.Range("T3").Value = -140
$ws.Range("U3").Value = 76
$ws.Range("Q4").Value = 0.002268304395076705
$ws.Range("R4").Value = 0.001143283054144875
$ws.Range("S4").Value = 0.002340601850973247
$ws.Range("T4").Value = 0.002931400170945582
$ws.Range("U4").Value = 0.004027559653477886
$ws.Range("Q6").Value = 4.3
$ws.Range("R6").Value = 4.4
$ws.Range("S6").Value = 4.5
$ws.Range("T6").ClearContents()
$ws.Range("U6").Value = 4.4
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8.4
$ws.Range("S7").Value = 8.7
$ws.Range("T7").ClearContents()
$ws.Range("U7").Value = 8.1
$ws.Range("Q8").Value = 62.5
$ws.Range("R8").Value = 62.4
$ws.Range("S8").Value = 62.5
$ws.Range("T8").ClearContents()
$ws.Range("U8").Value = 62.5
$ws.Range("Q9").Value = 59.8
$ws.Range("R9").Value = 59.7
$ws.Range("S9").Value = 59.6
$ws.Range("T9").ClearContents()
$ws.Range("U9").Value = 59.7
$ws.Range("Q15").Value = 34.3
$ws.Range("R15").Value = 34.2
$ws.Range("S15").Value = 34.3
$ws.Range("Q29").Value = 2.17
$ws.Range("R29").Value = 2.2
$ws.Range("S29").Value = 2.18
$ws.Range("T29").Value = 2.16
$ws.Range("Q30").Value = 2.32
$ws.Range("R30").Value = 2.35
$ws.Range("S30").Value = 2.34
$ws.Range("T30").Value = 2.32
$ws.Range("U30").Value = 2.35
$ws.Range("Q33").Value = 0.004051863857374327
$ws.Range("R33").Value = 0.000540540540540535
$ws.Range("S33").Value = 0.00407055630936215
$ws.Range("T33").Value = 0.004087193460490468
$ws.Range("U33").Value = 0.002184598580011077
$ws.Range("Q34").Value = 0.01128441879087681
$ws.Range("R34").Value = 0.0120947871629743
$ws.Range("S34").Value = 0.01032118383222275
$ws.Range("T34").Value = 0.01194571058798364
$ws.Range("U34").Value = 0.01322369465705905
$ws.Range("Q35").Value = 0.00199332572360067
$ws.Range("R35").Value = 0.002492631234120024
$ws.Range("S35").Value = -0.0004224555765300897
$ws.Range("T35").Value = 0.001486198108683112
$ws.Range("U35").Value = 0.001310931935978976
$ws.Range("Q36").Value = 0.01128441879087681
$ws.Range("R36").Value = 0.0120947871629743
$ws.Range("S36").Value = 0.01032118383222275
$ws.Range("T36").Value = 0.01194571058798364
$ws.Range("U36").Value = 0.01322369465705905
$ws.Range("Q48").Value = 3.48
$ws.Range("R48").Value = 3.5
$ws.Range("S48").Value = 3.47
$ws.Range("Q49").Value = 3.75
$ws.Range("R49").Value = 3.76
$ws.Range("S49").Value = 3.74
$ws.Range("Q50").Value = 4.22
$ws.Range("R50").Value = 4.22
$ws.Range("S50").Value = 4.21
$ws.Range("T50").Value = 4.29
$ws.Range("U50").Value = 4.28
$ws.Range("Q52").Value = 5.86
$ws.Range("R52").Value = 5.87
$ws.Range("S52").Value = 5.88
$ws.Range("T52").Value = 5.93
$ws.Range("U52").Value = 5.91
